# Updated symbol list on Tue Dec 27 06:53:54 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells in column D hold numeric-looking text that must stay text.
# Force the "Text" number format first so Excel does not silently convert
# the assigned string into a floating point number (with rounding drift).
$priceCells = @(
    "D2","D5","D7","D8","D9","D11","D12","D13","D14","D15","D16","D17",
    "D18","D19","D21","D23","D25","D40","D41","D42","D43","D44","D45",
    "D47","D48"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# --- Price (column D) updates ---
$ws.Range("D2").Value  = "243.42"
$ws.Range("D5").Value  = "0.05992"
$ws.Range("D7").Value  = "6.502"
$ws.Range("D8").Value  = "0.8091"
$ws.Range("D9").Value  = "0.9333"
$ws.Range("D11").Value = "0.07391"
$ws.Range("D12").Value = "0.03303"
$ws.Range("D13").Value = "0.03058"
$ws.Range("D14").Value = "0.09351"
$ws.Range("D15").Value = "3.859"
$ws.Range("D16").Value = "0.001579"
$ws.Range("D17").Value = "0.04696"
$ws.Range("D18").Value = "0.0005889"
$ws.Range("D19").Value = "0.005904"
$ws.Range("D21").Value = "0.004883"
$ws.Range("D23").Value = "3.574"
$ws.Range("D25").Value = "0.3238"
$ws.Range("D40").Value = "0.03958"

# --- Row 41: BKEXToken -> KickToken ---
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "0.006339"
$ws.Range("E41").Value = "40KickTokenKICK"

# --- Row 42: CEJI price/label update ---
$ws.Range("D42").Value = "0.004199"
$ws.Range("E42").Value = "41CEJICEJIBestin24h"

# --- Row 43: KickToken -> BKEXToken ---
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "0.1077"
$ws.Range("E43").Value = "42BKEXTokenBKK"

$ws.Range("D44").Value = "0.008540"
$ws.Range("D45").Value = "0.00005086"
$ws.Range("D47").Value = "0.6798"
$ws.Range("D48").Value = "0.002264"

# --- Column E label-only updates (no price changes in diff) ---
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E20").Value = "19BitKanKAN"
